# Remove the auto-inserted "Logo" pictures (Office sensitivity/compliance
# logo stamp) from every slide, and fix up the text on the "Overview and
# Project Goal" slide that used to reference "OCBang" by name.

$p = $ppt.ActivePresentation

# 1) Delete every picture shape whose alt text identifies it as the
#    "Logo ... Description automatically generated ..." stamp. Walk shapes
#    back-to-front so deleting doesn't perturb the indices still to visit.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Type -eq 13) {
            $alt = $shp.AlternativeText
            if ($alt -and $alt.StartsWith("Logo")) {
                $shp.Delete()
            }
        }
    }
}

# 2) Replace "technical recruitment company, OCBang Inc. also wants..."
#    with "technical recruitment company, we also wants..." on the
#    "Overview and Project Goal" slide, merging the three runs that spanned
#    the old text into a single run.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            $oldSub = "technical recruitment company, OCBang Inc. also wants to introduce this kind of system to fill the gap and improve the experience of clients. We aim to construct a bidirectional matching system between recruiter and potential candidates with machine learning techniques (especially advanced NLP techniques), improving the efficiency of recruitment activity and grabbing market share of our start-up."
            $pos = $full.IndexOf($oldSub)
            if ($pos -ge 0) {
                $newSub = "technical recruitment company, we also wants to introduce this kind of system to fill the gap and improve the experience of clients. We aim to construct a bidirectional matching system between recruiter and potential candidates with machine learning techniques (especially advanced NLP techniques), improving the efficiency of recruitment activity and grabbing market share of our start-up."
                $sub = $tr.Characters($pos + 1, $oldSub.Length)
                $sub.Text = $newSub
            }
        }
    }
}
